$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("2","0","Aircraft is returning to the Home Point Minimum RTH Altitude is 30m You can reset the RTH Altitude in Remote Controller Settings after cancelling RTH .","Minimum RTH Altitude is 30m","7-11","NonEvent","NonEvent"),
  @("3","0","Aircraft is returning to the Home Point Minimum RTH Altitude is 30m You can reset the RTH Altitude in Remote Controller Settings after cancelling RTH .","You can reset the RTH Altitude in Remote Controller Settings after cancelling RTH","12-24","NonEvent","NonEvent"),
  @("4","0","Aircraft is returning to the Home Point Minimum RTH Altitude is 30m You can reset the RTH Altitude in Remote Controller Settings after cancelling RTH .","Aircraft is returning to the Home Point","0-6","Event","Event"),
  @("5","1","Strong wind warning Aircraft unable to return to home automatically Lower altitude immediately and return to home manually .","Strong wind warning","0-2","Event","Event"),
  @("6","1","Strong wind warning Aircraft unable to return to home automatically Lower altitude immediately and return to home manually .","Lower altitude immediately and return to home manually","10-17","NonEvent","NonEvent"),
  @("7","1","Strong wind warning Aircraft unable to return to home automatically Lower altitude immediately and return to home manually .","Aircraft unable to return to home automatically","3-9","Event","Event"),
  @("8","2","Extra payload detected Max altitude and distance set Fly with caution .","Max altitude and distance set","3-7","Event","Event"),
  @("9","2","Extra payload detected Max altitude and distance set Fly with caution .","Extra payload detected","0-2","Event","Event"),
  @("10","2","Extra payload detected Max altitude and distance set Fly with caution .","Fly with caution","8-10","NonEvent","NonEvent"),
  @("11","3","Aircraft in Restricted Zone Unable to take off Check map to find Recommended Zones .","Check map to find Recommended Zones","8-13","NonEvent","NonEvent"),
  @("12","3","Aircraft in Restricted Zone Unable to take off Check map to find Recommended Zones .","Aircraft in Restricted Zone","0-3","Event","Event"),
  @("13","3","Aircraft in Restricted Zone Unable to take off Check map to find Recommended Zones .","Unable to take off","4-7","Event","Event"),
  @("14","4","Aircraft in Attitude mode Unable to hover Fly with caution .","Fly with caution","7-9","NonEvent","NonEvent"),
  @("15","4","Aircraft in Attitude mode Unable to hover Fly with caution .","Unable to hover","4-6","Event","Event"),
  @("16","4","Aircraft in Attitude mode Unable to hover Fly with caution .","Aircraft in Attitude mode","0-3","Event","Event"),
  @("17","5","GPS signal weak Positioning accuracy affected Fly with caution .","GPS signal weak","0-2","Event","Event"),
  @("18","5","GPS signal weak Positioning accuracy affected Fly with caution .","Fly with caution","6-8","NonEvent","NonEvent"),
  @("19","5","GPS signal weak Positioning accuracy affected Fly with caution .","Positioning accuracy affected","3-5","Event","Event"),
  @("20","6","Critically Low Power Aircraft Landing Throttle up to reduce the speed of descent and use sticks to avoid obstacles .","Critically Low Power","0-2","Event","Event"),
  @("21","6","Critically Low Power Aircraft Landing Throttle up to reduce the speed of descent and use sticks to avoid obstacles .","Aircraft Landing","3-4","Event","Event"),
  @("22","6","Critically Low Power Aircraft Landing Throttle up to reduce the speed of descent and use sticks to avoid obstacles .","Throttle up to reduce the speed of descent and use sticks to avoid obstacles","5-18","NonEvent","NonEvent"),
  @("23","7","Downward ambient light too low Obstacle avoidance unavailable Fly with caution .","Fly with caution","8-10","NonEvent","NonEvent"),
  @("24","7","Downward ambient light too low Obstacle avoidance unavailable Fly with caution .","Obstacle avoidance unavailable","5-7","NonEvent","NonEvent"),
  @("25","7","Downward ambient light too low Obstacle avoidance unavailable Fly with caution .","Downward ambient light too low","0-4","Event","Event"),
  @("26","8","Capture Failed Insufficient storage space Cannot shoot panorama .","Capture Failed","0-1","Event","Event"),
  @("27","8","Capture Failed Insufficient storage space Cannot shoot panorama .","Insufficient storage space","2-4","NonEvent","Event"),
  @("28","8","Capture Failed Insufficient storage space Cannot shoot panorama .","Cannot shoot panorama","5-7","Event","Event"),
  @("29","9","Approaching a No-Fly Zone RTH may be affected Fly with caution .","Fly with caution","8-10","NonEvent","NonEvent"),
  @("30","10","Remote ID error Connect to network to obtain mobile device location Flying without broadcasting required remote ID messages may violate local regulations .","Remote ID error","0-2","Event","Event"),
  @("31","10","Remote ID error Connect to network to obtain mobile device location Flying without broadcasting required remote ID messages may violate local regulations .","Flying without broadcasting required remote ID messages may violate local regulations","11-21","NonEvent","NonEvent"),
  @("32","10","Remote ID error Connect to network to obtain mobile device location Flying without broadcasting required remote ID messages may violate local regulations .","Connect to network to obtain mobile device location","3-10","NonEvent","NonEvent"),
  @("33","11","Compass Error Compass data error Please check the Compass installation and interference .","Compass Error","0-1","Event","Event"),
  @("34","11","Compass Error Compass data error Please check the Compass installation and interference .","Compass data error","2-4","Event","Event"),
  @("35","11","Compass Error Compass data error Please check the Compass installation and interference .","Please check the Compass installation and interference","5-11","NonEvent","NonEvent"),
  @("36","12","Weak signal Avoid blocking the antennas and keep the antennas parallel to and facing toward the aircraft during flight Downlink Lost .","Weak signal","0-1","Event","Event"),
  @("37","12","Weak signal Avoid blocking the antennas and keep the antennas parallel to and facing toward the aircraft during flight Downlink Lost .","Downlink Lost","19-20","Event","Event"),
  @("38","12","Weak signal Avoid blocking the antennas and keep the antennas parallel to and facing toward the aircraft during flight Downlink Lost .","Avoid blocking the antennas and keep the antennas parallel to and facing toward the aircraft during flight","2-18","NonEvent","NonEvent"),
  @("39","13","Aircraft is returning to the Home Point RTH Altitude is 30m You can reset the RTH Altitude in Remote Controller Settings after cancelling RTH if necessary .","RTH Altitude is 30m","7-10","NonEvent","NonEvent"),
  @("40","13","Aircraft is returning to the Home Point RTH Altitude is 30m You can reset the RTH Altitude in Remote Controller Settings after cancelling RTH if necessary .","You can reset the RTH Altitude in Remote Controller Settings after cancelling RTH if necessary","11-25","NonEvent","NonEvent"),
  @("41","13","Aircraft is returning to the Home Point RTH Altitude is 30m You can reset the RTH Altitude in Remote Controller Settings after cancelling RTH if necessary .","Aircraft is returning to the Home Point","0-6","Event","Event"),
  @("42","14","Battery Alert Battery installation error Please check the batteries are inserted correctly .","Please check the batteries are inserted correctly","5-11","NonEvent","NonEvent"),
  @("43","15","High altitude Aircraft braking distance increased and flight time decreased Fly with caution .","High altitude","0-1","NonEvent","NonEvent"),
  @("44","15","High altitude Aircraft braking distance increased and flight time decreased Fly with caution .","Aircraft braking distance increased and flight time decreased","2-9","Event","Event"),
  @("45","15","High altitude Aircraft braking distance increased and flight time decreased Fly with caution .","Fly with caution","10-12","NonEvent","NonEvent"),
  @("46","16","Aircraft entered Fixed-Wing Mode Current flight speed is nnn Press C2 to exit .","Aircraft entered Fixed-Wing Mode","0-3","Event","Event"),
  @("47","16","Aircraft entered Fixed-Wing Mode Current flight speed is nnn Press C2 to exit .","Current flight speed is nnn","4-8","NonEvent","Event"),
  @("48","16","Aircraft entered Fixed-Wing Mode Current flight speed is nnn Press C2 to exit .","Press C2 to exit","9-12","NonEvent","NonEvent"),
  @("49","17","Weak GPS signal Aircraft is in Attitude mode and hovering may be unstable Fly with caution .","Aircraft is in Attitude mode and hovering may be unstable","3-12","Event","NonEvent"),
  @("50","17","Weak GPS signal Aircraft is in Attitude mode and hovering may be unstable Fly with caution .","Weak GPS signal","0-2","Event","Event"),
  @("51","17","Weak GPS signal Aircraft is in Attitude mode and hovering may be unstable Fly with caution .","Fly with caution","13-15","NonEvent","NonEvent"),
  @("52","18","Camera error AI Spot-Check failed Restart camera .","Camera error","0-1","Event","Event"),
  @("53","18","Camera error AI Spot-Check failed Restart camera .","AI Spot-Check failed","2-4","Event","Event"),
  @("54","18","Camera error AI Spot-Check failed Restart camera .","Restart camera","5-6","NonEvent","NonEvent"),
  @("55","19","AirSense communication error Unable to display nearby manned aircraft and anti-collision warnings Fly with caution .","AirSense communication error","0-2","Event","Event"),
  @("56","19","AirSense communication error Unable to display nearby manned aircraft and anti-collision warnings Fly with caution .","Fly with caution","12-14","NonEvent","NonEvent"),
  @("57","19","AirSense communication error Unable to display nearby manned aircraft and anti-collision warnings Fly with caution .","Unable to display nearby manned aircraft and anti-collision warnings","3-11","Event","NonEvent"),
  @("58","20","Account not logged in Flight altitude and distance restricted to nnn and nnn Check and log in .","Check and log in","13-16","NonEvent","NonEvent")
)

foreach ($row in $data) {
    $r = [int]$row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
